# Updated ITA model - 2025-08-10 01:53
#
# The "Constants" sheet (xl/worksheets/sheet5.xml) gets a new discount-rate
# scenario ("USD21_alt") inserted into the G_DRATE table (columns G/I/J/K,
# rows 51-93). This shifts the existing G/I/J/K content down by one row
# (51->52 ... 93->94), turns row 51 into a lookup/formula row that mirrors
# the "base year" row (row 46) plus the literal "USD21_alt" label, and adds
# a brand-new trailing row 95 that echoes the new row 51 label back via a
# formula. Columns A/B (region labels) are untouched throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------
# 1) Shift the existing G/J/K content (old rows 51-93) down to rows 52-94.
#    Walk bottom-up so we never overwrite a source row before reading it.
# ---------------------------------------------------------------------
for ($r = 93; $r -ge 51; $r--) {
    $dest = $r + 1
    $ws.Range("G$dest").Value2 = $ws.Range("G$r").Value2
    $ws.Range("J$dest").Value2 = $ws.Range("J$r").Value2
    $ws.Range("K$dest").Value2 = $ws.Range("K$r").Value2
}

# ---------------------------------------------------------------------
# 2) Rebuild row 51 as the new "USD21_alt" entry, referencing the base
#    exchange-rate row (46) for the region/unit, same as the pattern used
#    by rows 40-50.
# ---------------------------------------------------------------------
$ws.Range("G51").Formula = "=G46"
$ws.Range("I51").Formula = "=I46"
$ws.Range("J51").Value2 = "USD21_alt"
$ws.Range("K51").Formula = "=K46"

# ---------------------------------------------------------------------
# 3) Add the brand-new trailing row 95, mirroring the new row 51 label.
# ---------------------------------------------------------------------
$ws.Range("G95").Value2 = "G_DRATE"
$ws.Range("J95").Formula = "=J51"
$ws.Range("K95").Value2 = 0.07

# ---------------------------------------------------------------------
# 4) View-state: the workbook now opens with the "Constants" tab active
#    (instead of "system_settings"), scrolled/selected on the new row.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("J95").Select()
